# Set line spacing to single (1x) for every paragraph in the document,
# matching Word's "Line Spacing = Single" setting
# (w:spacing w:line="240" w:lineRule="auto" in the underlying OOXML).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Format.LineSpacingRule = 0   # wdLineSpaceSingle
}
